$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest cryptos snapshot.
# Cells whose new value parses as a plain number need to be forced back to
# text (matching the original inlineStr cell type) and then have their style
# restored to Normal so no stray number formatting is left behind.

$ws.Range("D2").Value = "69.076.35"
$ws.Range("E2").Value = "  -3.68%  "
$ws.Range("D3").Value = "3.517.67"
$ws.Range("E3").Value = "  -4.63%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "3.509.16"
$ws.Range("E8").Value = "  -4.71%  "
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("E10").Value = "  -5.96%  "
$ws.Range("E11").Value = "  +7.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.601"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.09%  "
$ws.Range("E13").Value = "  -5.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000277"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "676.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.43%  "
$ws.Range("D16").Value = "4.084.77"
$ws.Range("E16").Value = "  -4.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.52%  "
$ws.Range("D18").Value = "3.525.37"
$ws.Range("E18").Value = "  -4.54%  "
$ws.Range("D19").Value = "69.068.02"
$ws.Range("E19").Value = "  -3.79%  "
$ws.Range("E20").Value = "  -1.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.909"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -9.05%  "
$ws.Range("E25").Value = "  -5.57%  "
$ws.Range("E26").Value = "  -4.17%  "
$ws.Range("E27").Value = "  -0.49%  "
$ws.Range("E28").Value = "  -6.14%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  -7.04%  "
$ws.Range("E31").Value = "  -7.06%  "
$ws.Range("E32").Value = "  -4.91%  "
$ws.Range("E33").Value = "  -8.02%  "
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("E35").Value = "  -5.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "580.31"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.08%  "
$ws.Range("E37").Value = "  -15.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.105"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "57.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.53%  "
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("E42").Value = "  -3.55%  "
$ws.Range("E43").Value = "  -4.69%  "
$ws.Range("E44").Value = "  -6.53%  "
$ws.Range("D45").Value = "3.434.58"
$ws.Range("E45").Value = "  -9.46%  "
$ws.Range("E46").Value = "  -5.30%  "
$ws.Range("D47").Value = "0.0₃0710"
$ws.Range("E47").Value = "  -8.99%  "
$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("E49").Value = "  -6.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "132.02"
$ws.Range("D51").Style = "Normal"
